$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.722.30"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.27%  '
$ws.Range("D3").Value = "'1.601.34"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.30%  '
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").Value = "'211.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.11%  '
$ws.Range("E6").Value = '  -0.38%  '
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.16%  '
$ws.Range("E8").Value = '  +0.19%  '
$ws.Range("E9").Value = '  +0.19%  '
$ws.Range("D10").Value = "'19.67"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.78%  '
$ws.Range("E11").Value = '  +0.91%  '
$ws.Range("D12").Value = "'1.825.87"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.26%  '
$ws.Range("D13").Value = "'1.607.42"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.06%  '
$ws.Range("E14").Value = '  +0.53%  '
$ws.Range("E15").Value = '  -0.30%  '
$ws.Range("D16").Value = "'65.09"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.04%  '
$ws.Range("E17").Value = '  +0.43%  '
$ws.Range("D18").Value = "'210.04"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.07%  '
$ws.Range("D19").Value = "'1.01"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.19%  '
$ws.Range("D20").Value = "'7.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.32%  '
$ws.Range("E21").Value = '  -0.24%  '
$ws.Range("E22").Value = '  -2.70%  '
$ws.Range("D23").Value = "'9.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.19%  '
$ws.Range("D24").Value = "'143.65"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.89%  '
$ws.Range("E25").Value = '  +0.10%  '
$ws.Range("E26").Value = '  -0.17%  '
$ws.Range("E27").Value = '  -0.83%  '
$ws.Range("D28").Value = "'15.38"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.69%  '
$ws.Range("D29").Value = "'0.0510"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.72%  '
$ws.Range("E30").Value = '  +0.09%  '
$ws.Range("E31").Value = '  +0.88%  '
$ws.Range("E32").Value = '  +0.91%  '
$ws.Range("D33").Value = "'1.292.10"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.80%  '
$ws.Range("E34").Value = '  +0.68%  '
$ws.Range("E35").Value = '  +0.60%  '
$ws.Range("D36").Value = "'0.603"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.58%  '
$ws.Range("D37").Value = "'1.17"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +10.37%  '
$ws.Range("E38").Value = '  -0.29%  '
$ws.Range("E39").Value = '  -0.37%  '
$ws.Range("E40").Value = '  -2.02%  '
$ws.Range("E41").Value = '  +0.12%  '
$ws.Range("D42").Value = "'0.785"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.09%  '
$ws.Range("D43").Value = "'62.69"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.65%  '
$ws.Range("D44").Value = "'1.737.96"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.21%  '
$ws.Range("D45").Value = "'90.56"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.06%  '
$ws.Range("E46").Value = '  -1.49%  '
$ws.Range("E47").Value = '  -0.33%  '
$ws.Range("D48").Value = "'0.0516"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.58%  '
$ws.Range("E49").Value = '  +0.12%  '
$ws.Range("D50").Value = "'7.43"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.86%  '
$ws.Range("E51").Value = '  +0.88%  '
